# Commit: "This is a test commit lol"
# Inserts a new bold / accent5-colored (with w14 shadow + text-outline
# text-effects) paragraph reading "OHHHH YEAH THESE ARE THE CHANGESSSSSS"
# right after the blank bordered paragraph that follows the "<Embed a
# copy of ... STS.>" instructions paragraph, and right before the final
# (empty) paragraph / sectPr.

$d = $word.ActiveDocument

# Locate the instructions paragraph by its known text, then the blank
# bordered paragraph immediately following it -- that's where the new
# paragraph belongs (just before the document's trailing empty <w:p/>).
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*submission of your STS*") {
        $anchorIndex = $i
    }
}

$blankBorderedPara = $d.Paragraphs.Item($anchorIndex + 1)

# Insert a brand-new paragraph right after it; this is the paragraph we
# will stamp with the new text + formatting.
$null = $blankBorderedPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 2)

# Build the run/paragraph-mark formatting exactly as authored: bold,
# accent5 themed blue, plus the w14 text-effect extensions (soft shadow
# and a thin white text outline) that Word stamps on both the paragraph
# mark's rPr and the run's rPr.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:rPr><w:b/><w:color w:val="5B9BD5" w:themeColor="accent5"/><w14:shadow w14:blurRad="12700" w14:dist="38100" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl"><w14:schemeClr w14:val="accent5"><w14:lumMod w14:val="60000"/><w14:lumOff w14:val="40000"/></w14:schemeClr></w14:shadow><w14:textOutline w14:w="9525" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="bg1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="5B9BD5" w:themeColor="accent5"/><w14:shadow w14:blurRad="12700" w14:dist="38100" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl"><w14:schemeClr w14:val="accent5"><w14:lumMod w14:val="60000"/><w14:lumOff w14:val="40000"/></w14:schemeClr></w14:shadow><w14:textOutline w14:w="9525" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="bg1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr><w:t>OHHHH YEAH THESE ARE THE CHANGESSSSSS</w:t></w:r></w:p>'

$null = $newPara.Range.InsertXML($xml)

Write-Output "Inserted new paragraph at index $($anchorIndex + 2)"
